$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 26266.205
$ws.Range("I15").Value = 26266.205
$ws.Range("K15").Value = 78798.61500000001
$ws.Range("M15").Value = -78629.61500000001
$ws.Range("H40").Value = 27484.75
$ws.Range("J40").Value = 3474.5
$ws.Range("L40").Value = 3474.5
$ws.Range("N40").Value = -3824.5
$ws.Range("H42").Value = 198
$ws.Range("I42").Value = 151
$ws.Range("J42").Value = 254.4
$ws.Range("K42").Value = 453
$ws.Range("L42").Value = 763.2
$ws.Range("M42").Value = -223
$ws.Range("N42").Value = -1223.2
$ws.Range("H129").Value = 1478.3077
$ws.Range("J129").Value = 2403.5715
$ws.Range("L129").Value = 7210.7145
$ws.Range("N129").Value = -17210.7145
$ws.Range("H137").Value = 2308
$ws.Range("I137").Value = 2148.111
$ws.Range("J137").Value = 2392.647
$ws.Range("K137").Value = 6444.333
$ws.Range("L137").Value = 7177.941
$ws.Range("M137").Value = -3894.333
$ws.Range("N137").Value = -12277.941
$ws.Range("H138").Value = 4170954.5
$ws.Range("J138").Value = 7148261
$ws.Range("L138").Value = 21444783
$ws.Range("N138").Value = -21455063

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7327.8423
$ws.Range("I45").Value = 2737.75
$ws.Range("K45").Value = 2737.75
$ws.Range("M45").Value = -2360.75
$ws.Range("H74").Value = 60895.945
$ws.Range("I74").Value = 113714.555
$ws.Range("K74").Value = 113714.555
$ws.Range("M74").Value = -112840.555
$ws.Range("H77").Value = 60895.945
$ws.Range("I77").Value = 113714.555
$ws.Range("K77").Value = 568572.7749999999
$ws.Range("M77").Value = -564204.7749999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 7966.769
$ws.Range("J22").Value = 325
$ws.Range("L22").Value = 325
$ws.Range("N22").Value = -671
$ws.Range("H110").Value = 59384
$ws.Range("J110").Value = 59384
$ws.Range("L110").Value = 59384
$ws.Range("N110").Value = -67564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3639.7551
$ws.Range("I132").Value = 1687.0286
$ws.Range("K132").Value = 5061.085800000001
$ws.Range("M132").Value = -2531.085800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 142857920
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 142857920
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 428573760
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -428574326
$ws.Range("H51").Value = 2800
$ws.Range("I51").Value = 100
$ws.Range("J51").Value = 3700
$ws.Range("K51").Value = 300
$ws.Range("L51").Value = 11100
$ws.Range("M51").Value = 160
$ws.Range("N51").Value = -12020
$ws.Range("H52").Value = 1595.4286
$ws.Range("J52").Value = 1595.4286
$ws.Range("L52").Value = 4786.2858
$ws.Range("N52").Value = -5318.2858
$ws.Range("H131").Value = 2054.4
$ws.Range("I131").Value = 1171.4445
$ws.Range("J131").Value = 2275.139
$ws.Range("K131").Value = 3514.3335
$ws.Range("L131").Value = 6825.417
$ws.Range("M131").Value = 1525.6665
$ws.Range("N131").Value = -16905.417
$ws.Range("H141").Value = 2836.4443
$ws.Range("I141").Value = 2836.4443
$ws.Range("K141").Value = 8509.332900000001
$ws.Range("M141").Value = -3329.332900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 14934.25
$ws.Range("J36").Value = 15868.5
$ws.Range("L36").Value = 15868.5
$ws.Range("N36").Value = -16838.5
$ws.Range("H64").Value = 70000
$ws.Range("I64").Value = 40000
$ws.Range("J64").Value = 100000
$ws.Range("K64").Value = 40000
$ws.Range("L64").Value = 100000
$ws.Range("M64").Value = -39752
$ws.Range("N64").Value = -100496
$ws.Range("H67").Value = 70000
$ws.Range("I67").Value = 40000
$ws.Range("J67").Value = 100000
$ws.Range("K67").Value = 40000
$ws.Range("L67").Value = 100000
$ws.Range("M67").Value = -39142
$ws.Range("N67").Value = -101716
$ws.Range("H80").Value = 3130.7778
$ws.Range("I80").Value = 1966.6
$ws.Range("J80").Value = 4586
$ws.Range("K80").Value = 1966.6
$ws.Range("L80").Value = 4586
$ws.Range("M80").Value = -968.5999999999999
$ws.Range("N80").Value = -6582
$ws.Range("H83").Value = 3130.7778
$ws.Range("I83").Value = 1966.6
$ws.Range("J83").Value = 4586
$ws.Range("K83").Value = 9833
$ws.Range("L83").Value = 22930
$ws.Range("M83").Value = -4841
$ws.Range("N83").Value = -32914
$ws.Range("H117").Value = 56961
$ws.Range("J117").Value = 56961
$ws.Range("L117").Value = 56961
$ws.Range("N117").Value = -63845
$ws.Range("H126").Value = 5330.3076
$ws.Range("I126").Value = 2573.5
$ws.Range("K126").Value = 7720.5
$ws.Range("M126").Value = -5250.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2420.0908
$ws.Range("I16").Value = 2532.3157
$ws.Range("J16").Value = 1709.3334
$ws.Range("K16").Value = 2532.3157
$ws.Range("L16").Value = 1709.3334
$ws.Range("M16").Value = -2362.3157
$ws.Range("N16").Value = -2049.3334
$ws.Range("H22").Value = 1339.238
$ws.Range("I22").Value = 1065.3572
$ws.Range("J22").Value = 1887
$ws.Range("K22").Value = 1065.3572
$ws.Range("L22").Value = 1887
$ws.Range("M22").Value = -770.3571999999999
$ws.Range("N22").Value = -2477
$ws.Range("H27").Value = 1339.238
$ws.Range("I27").Value = 1065.3572
$ws.Range("J27").Value = 1887
$ws.Range("K27").Value = 1065.3572
$ws.Range("L27").Value = 1887
$ws.Range("M27").Value = -958.3571999999999
$ws.Range("N27").Value = -2101
$ws.Range("H115").Value = 57976
$ws.Range("J115").Value = 57976
$ws.Range("L115").Value = 57976
$ws.Range("N115").Value = -60326
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 6229.587
$ws.Range("I122").Value = 6111.9355
$ws.Range("J122").Value = 6472.7334
$ws.Range("K122").Value = 18335.8065
$ws.Range("L122").Value = 19418.2002
$ws.Range("M122").Value = -15885.8065
$ws.Range("N122").Value = -24318.2002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3354.2727
$ws.Range("I96").Value = 2899.625
$ws.Range("K96").Value = 2899.625
$ws.Range("M96").Value = -1526.625
$ws.Range("H116").Value = 56965
$ws.Range("J116").Value = 56965
$ws.Range("L116").Value = 56965
$ws.Range("N116").Value = -66143
$ws.Range("H132").Value = 16080.535
$ws.Range("I132").Value = 14185.64
$ws.Range("K132").Value = 42556.92
$ws.Range("M132").Value = -40026.92
